$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "test"

# Fix typo in header, insert a new "type" column between price and user,
# and split the old combined "user" value into "type" + "user".
$ws.Columns.Item(5).Insert()

$ws.Range("C1").Value = "weight"
$ws.Range("E1").Value = "type"
$ws.Range("F1").Value = "user"

$ws.Range("E2").Value = "S"
$ws.Range("F2").Value = "u1"

$ws.Range("E3").Value = "W"
$ws.Range("F3").Value = "u2"

# Column widths as specified in the target layout (nearest reachable value
# given the host's pixel-grid rounding of ColumnWidth -> stored XML width)
$ws.Columns.Item(1).ColumnWidth = 5.6667
$ws.Columns.Item(2).ColumnWidth = 7.1667
$ws.Columns.Item(3).ColumnWidth = 27.6667
$ws.Columns.Item(4).ColumnWidth = 33
$ws.Columns.Item(5).ColumnWidth = 4.1667
$ws.Columns.Item(6).ColumnWidth = 4.1667

$ws.Range("I10").Select()
